$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B for the two new rating dates
# (Jun_26 appears twice, Jun_27 once) - shifts old B:E (Jun_17..Jun_10) to E:H
$ws.Columns("B:D").Insert()

# Keep the same custom column width (8.0) across the now-wider date block
$ws.Columns("C:G").ColumnWidth = 7.1

$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"
$ws.Range("B2:D27").Value = "UN"

# New analyst group rows
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28:D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29:D29").Value = "UN"
